$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: dimension/measure metadata type changes (curated dimensions)
$ws.Range("C2").Value = "iaest-measure:tipo-personal"
$ws.Range("I2").Value = "iaest-measure:jornadas-trabajo"
$ws.Range("J2").Value = "sdmx-dimension:refArea"
$ws.Range("L2").Value = "iaest-measure:sexo"

# Row 3: dim/medida role changes
$ws.Range("C3").Value = "medida"
$ws.Range("I3").Value = "medida"
$ws.Range("J3").Value = "dim"
$ws.Range("L3").Value = "medida"

# Row 4: data type / URI changes
$ws.Range("C4").Value = "xsd:int"
$ws.Range("I4").Value = "xsd:int"
$ws.Range("J4").Value = "URI-Municipio"
$ws.Range("L4").Value = "xsd:int"

# Row 5: remove obsolete mapping file references entirely
$ws.Range("C5").Clear()
$ws.Range("I5").Clear()
$ws.Range("L5").Clear()
